$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all target cells are treated as plain Text so Excel does not
# auto-convert numeric-looking strings (e.g. "1.00", "0.0000237") into
# floating point numbers, which would lose the original text formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.272.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.056.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.18"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.053.66"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.90"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.555.32"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.203.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.055.69"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.96"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.705"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.51"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.88"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.82"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.50"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.68"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.23"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.79%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0819"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.03"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.31"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.46"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "438.27"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.81%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0362"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.820.86"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.43"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.03"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.36"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.30%  "
